# Apply the fitting-parameter update + active-sheet/selection changes.
$wb = $excel.ActiveWorkbook

# --- Parameters sheet: update fitted values in J2 (r_s_star) and K2 (h_p_star) ---
$wsParams = $wb.Worksheets.Item("Parameters")
$wsParams.Range("J2").Value = 0.0337
$wsParams.Range("K2").Value = 0.23866

# --- Selection / active sheet bookkeeping ---
# Parameters becomes the active / tab-selected sheet, with a new selection anchor.
# (DataExp keeps its own prior selection untouched; it simply stops being the
# active tab once another sheet is selected.)
$wsParams.Select()
$wsParams.Range("K3").Select()
